# Trading update: 2026-02-17 13:24:19
# Append a new OPEN trade (row 40) to the "All Trades" and "MarketMaking"
# sheets. Both sheets carry the same trade log, so the new row is written
# identically to each one.

$wb = $excel.ActiveWorkbook

# ---- "All Trades" sheet ----
$ws = $wb.Worksheets.Item("All Trades")

# A: Trade #
$ws.Cells.Item(40, 1).Value = 39

# B: Date - write as literal text (quote-prefixed so the ISO-looking date
# string is not auto-converted into a date serial number), then restore
# the Normal style so no stray number format is left behind on the cell.
$ws.Cells.Item(40, 2).NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = "2026-02-17"
$ws.Cells.Item(40, 2).Style = "Normal"

# C: Time (plain text; no auto-conversion risk)
$ws.Cells.Item(40, 3).Value = "13:23:45"

# D: Strategy
$ws.Cells.Item(40, 4).Value = "MarketMaking"

# E: Side
$ws.Cells.Item(40, 5).Value = "DOWN"

# F: Entry Price
$ws.Cells.Item(40, 6).Value = 0.98

# G: Exit Price - left blank; trade is still OPEN, no exit price yet.

# H: Status
$ws.Cells.Item(40, 8).Value = "OPEN"

# I: P&L %
$ws.Cells.Item(40, 9).Value = 0

# J: P&L $
$ws.Cells.Item(40, 10).Value = 0

# K: Capital After
$ws.Cells.Item(40, 11).Value = 98.32811449838626

# L: Entry Slippage (bps)
$ws.Cells.Item(40, 12).Value = 0

# M: Exit Slippage (bps)
$ws.Cells.Item(40, 13).Value = 0

# N: Confidence
$ws.Cells.Item(40, 14).Value = 0.6

# O: Entry Reason
$ws.Cells.Item(40, 15).Value = "Normal spread capture: 19600 bps"

# P: Exit Reason - left blank; trade is still OPEN.

# Q: Duration (min)
$ws.Cells.Item(40, 17).Value = 0

# ---- "MarketMaking" sheet (identical new row) ----
$ws = $wb.Worksheets.Item("MarketMaking")

$ws.Cells.Item(40, 1).Value = 39

$ws.Cells.Item(40, 2).NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = "2026-02-17"
$ws.Cells.Item(40, 2).Style = "Normal"

$ws.Cells.Item(40, 3).Value = "13:23:45"
$ws.Cells.Item(40, 4).Value = "MarketMaking"
$ws.Cells.Item(40, 5).Value = "DOWN"
$ws.Cells.Item(40, 6).Value = 0.98

$ws.Cells.Item(40, 8).Value = "OPEN"
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 98.32811449838626
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = 0
$ws.Cells.Item(40, 14).Value = 0.6
$ws.Cells.Item(40, 15).Value = "Normal spread capture: 19600 bps"

$ws.Cells.Item(40, 17).Value = 0
